$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$headers = @(
    "rollout",
    "waitingTime",
    "time",
    "rew",
    "E_num_train_rollouts",
    "E_rollout_length",
    "E_eval_freq",
    "E_eval_num_eps",
    "E_max_ep_steps",
    "E_generation_ep_steps",
    "E_test_num_eps",
    "P_gae_tau",
    "P_entropy_weight",
    "P_minibatch_size",
    "P_optimization_epochs",
    "P_ppo_ratio_clip",
    "P_discount",
    "P_learning_rate",
    "P_clip_grads",
    "P_gradient_clip",
    "P_value_loss_coef",
    "R_rule_set",
    "R_rule_set_params",
    "O_num_workers",
    "O_environment",
    "O_agent_type",
    "O_single_agent"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Add new data row (row 2) ---
$rowValues = @(
    250,
    1769.131666666667,
    31.3077213605245,
    426.59525,
    400,
    2048,
    50,
    25,
    900,
    750,
    50,
    0.85,
    0.1,
    256,
    5,
    0.3,
    0.99,
    0.0001,
    $true,
    1,
    1,
    "cycle",
    "{'cycle_length': 20, 'NS_mult': 1.5, 'EW_mult': 1.0, 'phase_end_offset': 100}",
    8,
    "4_intersections_rush_hour",
    "ppo",
    $false
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $rowValues[$i]
}
